$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.430.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.522.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.518.75"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.05%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000222"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.114.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.505.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.399.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.34%  "
$ws.Range("E20").Value = "  -1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "454.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.645"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.660.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.06%  "
$ws.Range("E27").Value = "  -4.43%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.34%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.157"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.518.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.40%  "
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.889"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.78%  "
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.86%  "
